# Minor changes to getting started guide
#
# 1) Refresh the cached "datetimeFigureOut" date field text (slide master +
#    all 11 slide layouts) from 8/12/2019 -> 9/16/2019.
# 2) On slide 4, rename the $pack-x-y / $wh-x placeholder tokens to use
#    underscores instead of hyphens, and nudge a few of the textbox
#    positions/sizes to match their (slightly) widened labels.

$p = $ppt.ActivePresentation

# --- 1) Date placeholders on the master + every custom layout ---------

$EMU_PER_PT = 12700.0
$newDate = "9/16/2019"

function Set-DatePlaceholderText($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Reach the master through Designs (Slides.Item(n).Master always resolves
# to the layout used by the "current" slide in this host, so go through
# the Design -> SlideMaster path instead to properly enumerate all of the
# distinct custom layouts).
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Set-DatePlaceholderText $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes
}

# --- 2) Slide 4 product-pack / warehouse textboxes ---------------------

$s4 = $p.Slides.Item(4)

function Set-ShapeFrame($shape, $x, $y, $cx, $cy) {
    $shape.Left = $x / $EMU_PER_PT
    $shape.Top = $y / $EMU_PER_PT
    $shape.Width = $cx / $EMU_PER_PT
    $shape.Height = $cy / $EMU_PER_PT
}

for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shp = $s4.Shapes.Item($i)
    if (-not $shp.HasTextFrame) {
        continue
    }
    $txt = $shp.TextFrame.TextRange.Text

    if ($txt -eq "`$pack-1-1") {
        Set-ShapeFrame $shp 6536501 3138850 540704 369332
        $shp.TextFrame.TextRange.Text = "`$pack_1_1"
    }
    elseif ($txt -eq "`$pack-1-3") {
        Set-ShapeFrame $shp 7103275 3138850 537557 369332
        $shp.TextFrame.TextRange.Text = "`$pack_1_3"
    }
    elseif ($txt -eq "`$wh-1   `$wh-2   `$wh-3") {
        $shp.TextFrame.TextRange.Text = "`$wh_1   `$wh_2   `$wh_3"
    }
    elseif ($txt -eq "`$pack-2-1") {
        Set-ShapeFrame $shp 6536501 759994 537557 369332
        $shp.TextFrame.TextRange.Text = "`$pack_2_1"
    }
    elseif ($txt -eq "`$pack-2-3") {
        Set-ShapeFrame $shp 7074058 759994 537557 369332
        $shp.TextFrame.TextRange.Text = "`$pack_2_3"
    }
}
